# Update marksheet marks: correct-answer marking scheme changed from 3 to 5,
# which cascades into the "Total" row marks and the "Correct/Total" fraction.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Marking (per correct answer) : 3 -> 5
$ws.Range("B11").Value = 5

# Total marks obtained : 63 -> 105
$ws.Range("B12").Value = 105

# Correct/Total marks text : "61/84" -> "105/140"
$ws.Range("E12").Value = "105/140"
